$wb = $excel.ActiveWorkbook

# --- Sheet "Full results" ---
$ws1 = $wb.Worksheets.Item("Full results")

$ws1.Range("C2").Value = 0.675903707073575
$ws1.Range("D2").Value = 0.324438889402782
$ws1.Range("E2").Value = 1.00034259647636
$ws1.Range("J2").Value = 0.324327775849591
$ws1.Range("K2").Value = 0.287203319078383
$ws1.Range("L2").Value = 0.00933406479094052
$ws1.Range("M2").Value = 0.034262491581737
$ws1.Range("N2").Value = 0.296537383869324

$ws1.Range("F3").Value = 0.650966739891599
$ws1.Range("G3").Value = 0.287301713923497

$ws1.Range("H4").Value = 0.641629477282951
$ws1.Range("I4").Value = 0.241756267272275
$ws1.Range("O4").Value = 0.358590267431328

# --- Sheet "For plotting" ---
$ws2 = $wb.Worksheets.Item("For plotting")

$ws2.Range("C2").Value = 0.324327775849591
$ws2.Range("D2").Value = 0.264162236920513
$ws2.Range("E2").Value = 0.38449331477867
$ws2.Range("F2").Value = 948

$ws2.Range("C3").Value = 0.296537383869324
$ws2.Range("D3").Value = 0.217012782397054
$ws2.Range("E3").Value = 0.376061985341593
$ws2.Range("F3").Value = 948

$ws2.Range("C4").Value = 0.358590267431328
$ws2.Range("D4").Value = 0.279364345887865
$ws2.Range("E4").Value = 0.437816188974792
$ws2.Range("F4").Value = 948
